$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "ValidLogin"

# Populate header and test data
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

# Auto-fit column A width to match content (bestFit)
$ws.Columns.Item(1).AutoFit() | Out-Null

# Update selection to B2 to match the saved view state
$ws.Range("B2").Select()
